$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F13").Value = 351
$wsExhibit.Range("F15").Value = 3080
$wsExhibit.Range("F18").Value = 1728

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F14").Value = 351
$wsAll.Range("F16").Value = 3080
$wsAll.Range("F19").Value = 1728
